$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number must be forced to stay text
# (matching the original inlineStr text cells), then restored to the default "Normal"
# style so no stray number formatting is left behind.

$ws.Range("D2").Value = "26.676.44"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.634.91"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").Value = "1.863.69"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "1.639.14"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "26.662.10"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.29%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("E28").Value = "  +5.74%  "
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").Value = "1.199.63"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E37").Value = "  +5.73%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.793"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "1.771.13"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  +0.04%  "
